$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column X: "Challenge_class" -------------------------------------
# Header cell (bold font, wrapped, matching the header row style used by
# the rest of row 1) plus the per-row challenge-class numbers for rows 2-70.

$ws.Range("X1").Value = "Challenge_class"
$ws.Range("X1").Font.Bold = $true
$ws.Range("X1").Font.Color = 65536
$ws.Range("X1").WrapText = $true

$classValues = @(1,1,2,2,3,3,4,5,6,7,8,9,9,9,10,10,11,11,11,12,13,14,15,15,16,17,18,18,19,19,19,20,20,20,21,22,23,24,24,24,24,25,26,27,27,27,27,28,29,30,30,30,31,31,32,33,34,35,35,36,36,36,37,38,39,40,40,41,42)

for ($i = 0; $i -lt $classValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 24).Value = $classValues[$i]
}

# --- Column width for the new column (best effort in this environment) ---
$ws.Columns.Item(24).ColumnWidth = 15.92
